$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the date values between row 3 and row 4 (column D)
$ws.Range("D3").Value = 44414
$ws.Range("D4").Value = 44379

# Swap the volume values between row 3 and row 4 (column J)
$ws.Range("J3").Value = 500
$ws.Range("J4").Value = 240
